$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update transaction amount column (B) ---
$ws.Range("B3").Value = 1.88
$ws.Range("B4").Value = 1.88
$ws.Range("B5").Value = 1.88

# --- Update ChainTag column (C) ---
$ws.Range("C3").Value = "0x27"
$ws.Range("C4").Value = "0x27"
$ws.Range("C5").Value = "0x27"

# --- Update BlockRef column (D) ---
$ws.Range("D3").Value = "0x0000695540f491a5"
$ws.Range("D4").Value = "0x0000695540f491a5"
$ws.Range("D5").Value = "0x0000695540f491a5"

# D column no longer uses the special BlockRef font -- align it back to the
# plain font used by the rest of the data rows (matches column C).
$ws.Range("D3:D5").Font.Name = "Helvetica Neue"
$ws.Range("D3:D5").Font.Size = 10

# --- Update the "To" address on the third data row ---
$ws.Range("A5").Value = "0xf881a94423f22ee9a0e3e1442f515f43c966b7ed"

# --- Move the active selection like the refreshed workbook does ---
[void]$ws.Range("A8").Select()
